$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume/name/link data per scheduled GitHub Actions refresh.
# NumberFormat is forced to Text ("@") before assigning so that numeric-looking
# strings (e.g. "585.30") are stored as text, matching the source data feed,
# then the cell Style is reset to "Normal" so no stray custom format is left behind.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "63.907.53"
Set-TextValue 2 5 "  -1.85%  "

Set-TextValue 3 4 "3.507.93"
Set-TextValue 3 5 "  -0.49%  "

Set-TextValue 4 5 "  +0.13%  "

Set-TextValue 5 4 "585.30"
Set-TextValue 5 5 "  -1.31%  "

Set-TextValue 6 4 "132.89"
Set-TextValue 6 5 "  -1.20%  "

Set-TextValue 7 4 "3.512.13"
Set-TextValue 7 5 "  -0.35%  "

Set-TextValue 8 5 "  +0.01%  "

Set-TextValue 9 4 "0.485"
Set-TextValue 9 5 "  -1.19%  "

Set-TextValue 10 5 "  -0.87%  "

Set-TextValue 11 4 "7.10"
Set-TextValue 11 5 "  -0.71%  "

Set-TextValue 12 4 "0.376"
Set-TextValue 12 5 "  -2.65%  "

Set-TextValue 13 4 "4.118.51"
Set-TextValue 13 5 "  -0.12%  "

Set-TextValue 14 4 "27.43"
Set-TextValue 14 5 "  -0.63%  "

Set-TextValue 15 5 "  +1.21%  "

Set-TextValue 16 4 "3.533.79"
Set-TextValue 16 5 "  +0.16%  "

Set-TextValue 17 4 "0.0000178"
Set-TextValue 17 5 "  -1.93%  "

Set-TextValue 18 4 "64.055.96"
Set-TextValue 18 5 "  -1.54%  "

Set-TextValue 19 4 "9.78"
Set-TextValue 19 5 "  -2.85%  "

Set-TextValue 20 4 "13.92"
Set-TextValue 20 5 "  -3.31%  "

Set-TextValue 21 4 "5.60"
Set-TextValue 21 5 "  -1.34%  "

Set-TextValue 22 4 "383.93"
Set-TextValue 22 5 "  -2.11%  "

Set-TextValue 23 4 "0.571"
Set-TextValue 23 5 "  -1.51%  "

Set-TextValue 24 4 "3.657.98"
Set-TextValue 24 5 "  -0.22%  "

Set-TextValue 25 4 "73.85"
Set-TextValue 25 5 "  -1.19%  "

Set-TextValue 26 5 "  +0.06%  "

Set-TextValue 27 4 "5.66"
Set-TextValue 27 5 "  -0.46%  "

Set-TextValue 28 5 "  +1.88%  "

Set-TextValue 29 4 "1.56"
Set-TextValue 29 5 "  -1.53%  "

Set-TextValue 30 4 "7.47"
Set-TextValue 30 5 "  -2.89%  "

Set-TextValue 31 5 "  +0.40%  "

Set-TextValue 32 4 "8.37"
Set-TextValue 32 5 "  +0.38%  "

Set-TextValue 33 4 "2.22"
Set-TextValue 33 5 "  -2.20%  "

Set-TextValue 34 4 "3.526.59"
Set-TextValue 34 5 "  -0.12%  "

Set-TextValue 36 4 "23.49"
Set-TextValue 36 5 "  -2.61%  "

Set-TextValue 37 5 "  -0.09%  "

Set-TextValue 38 4 "5.33"
Set-TextValue 38 5 "  +0.40%  "

Set-TextValue 39 5 "  -0.75%  "

Set-TextValue 40 4 "6.89"
Set-TextValue 40 5 "  -0.79%  "

Set-TextValue 41 4 "160.74"
Set-TextValue 41 5 "  -4.45%  "

Set-TextValue 42 4 "0.0783"
Set-TextValue 42 5 "  -2.91%  "

Set-TextValue 43 2 "EnergySwap"
Set-TextValue 43 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue 43 4 "26.67"
Set-TextValue 43 5 "  +3.94%  "

Set-TextValue 44 2 "Mantle"
Set-TextValue 44 3 "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue 44 4 "0.810"
Set-TextValue 44 5 "  -1.58%  "

Set-TextValue 45 5 "  +0.11%  "

Set-TextValue 46 2 "OKB"
Set-TextValue 46 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue 46 4 "41.55"
Set-TextValue 46 5 "  -3.36%  "

Set-TextValue 47 2 "ONDO"
Set-TextValue 47 3 "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue 47 4 "1.21"
Set-TextValue 47 5 "  -4.44%  "

Set-TextValue 48 4 "4.39"
Set-TextValue 48 5 "  -1.05%  "

Set-TextValue 49 4 "1.61"
Set-TextValue 49 5 "  -3.74%  "

Set-TextValue 50 4 "2.481.24"
Set-TextValue 50 5 "  +2.44%  "

Set-TextValue 51 4 "6.78"
Set-TextValue 51 5 "  -1.77%  "
